# Refresh the crypto price/volume snapshot (cryptos.xlsx), per the
# "Updated cryptos list ... with GitHub Actions" scrape commit.
#
# Column D (Price) / E (Volume 1h) cells are stored as literal text in the
# workbook (several prices even use "." as a thousands separator, e.g.
# "26.302.01"), so a handful of the new D-column values look like plain
# decimals ("1.000", "18.56", ...) that Excel's normal Range.Value setter
# would silently reinterpret as numbers (dropping trailing zeros / exponent
# notation). For those specific values we route the write through a
# throwaway cell using ="<text>" + Copy/PasteSpecial(values-only), which
# lands the literal string without leaving any NumberFormat/style residue
# behind (unlike forcing NumberFormat = "@" on the destination cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$scratch = "ZZ1"   # far outside the A1:E51 table, cleared after each use


# Row 2
$ws.Range("D2").Value = "26.302.01"
$ws.Range("E2").Value = "  -3.06%  "

# Row 3
$ws.Range("D3").Value = "1.830.75"
$ws.Range("E3").Value = "  -2.73%  "

# Row 4
$ws.Range("E4").Value = "  +0.10%  "

# Row 5
$ws.Range($scratch).Formula = "=""258.31"""
$ws.Range($scratch).Copy()
$ws.Range("D5").PasteSpecial(-4163)  # xlPasteValues
$ws.Range($scratch).Clear()
$ws.Range("E5").Value = "  -8.34%  "

# Row 6
$ws.Range($scratch).Formula = "=""1.000"""
$ws.Range($scratch).Copy()
$ws.Range("D6").PasteSpecial(-4163)  # xlPasteValues
$ws.Range($scratch).Clear()
$ws.Range("E6").Value = "  +0.07%  "

# Row 7
$ws.Range($scratch).Formula = "=""0.5178"""
$ws.Range($scratch).Copy()
$ws.Range("D7").PasteSpecial(-4163)  # xlPasteValues
$ws.Range($scratch).Clear()
$ws.Range("E7").Value = "  -2.20%  "

# Row 8
$ws.Range($scratch).Formula = "=""0.3218"""
$ws.Range($scratch).Copy()
$ws.Range("D8").PasteSpecial(-4163)  # xlPasteValues
$ws.Range($scratch).Clear()
$ws.Range("E8").Value = "  -9.10%  "

# Row 9
$ws.Range("E9").Value = "  -4.77%  "

# Row 10
$ws.Range($scratch).Formula = "=""18.56"""
$ws.Range($scratch).Copy()
$ws.Range("D10").PasteSpecial(-4163)  # xlPasteValues
$ws.Range($scratch).Clear()
$ws.Range("E10").Value = "  -8.95%  "

# Row 11
$ws.Range($scratch).Formula = "=""0.7628"""
$ws.Range($scratch).Copy()
$ws.Range("D11").PasteSpecial(-4163)  # xlPasteValues
$ws.Range($scratch).Clear()
$ws.Range("E11").Value = "  -7.25%  "

# Row 12
$ws.Range($scratch).Formula = "=""0.07670"""
$ws.Range($scratch).Copy()
$ws.Range("D12").PasteSpecial(-4163)  # xlPasteValues
$ws.Range($scratch).Clear()
$ws.Range("E12").Value = "  -1.82%  "

# Row 13
$ws.Range("D13").Value = "1.868.12"
$ws.Range("E13").Value = "  -0.69%  "

# Row 14
$ws.Range($scratch).Formula = "=""88.43"""
$ws.Range($scratch).Copy()
$ws.Range("D14").PasteSpecial(-4163)  # xlPasteValues
$ws.Range($scratch).Clear()
$ws.Range("E14").Value = "  -2.58%  "

# Row 15
$ws.Range($scratch).Formula = "=""5.010"""
$ws.Range($scratch).Copy()
$ws.Range("D15").PasteSpecial(-4163)  # xlPasteValues
$ws.Range($scratch).Clear()
$ws.Range("E15").Value = "  -3.93%  "

# Row 16
$ws.Range("E16").Value = "  +0.11%  "

# Row 17
$ws.Range($scratch).Formula = "=""14.02"""
$ws.Range($scratch).Copy()
$ws.Range("D17").PasteSpecial(-4163)  # xlPasteValues
$ws.Range($scratch).Clear()
$ws.Range("E17").Value = "  -4.16%  "

# Row 18
$ws.Range($scratch).Formula = "=""1.000"""
$ws.Range($scratch).Copy()
$ws.Range("D18").PasteSpecial(-4163)  # xlPasteValues
$ws.Range($scratch).Clear()
$ws.Range("E18").Value = "  +0.08%  "

# Row 19
$ws.Range($scratch).Formula = "=""0.000007879"""
$ws.Range($scratch).Copy()
$ws.Range("D19").PasteSpecial(-4163)  # xlPasteValues
$ws.Range($scratch).Clear()
$ws.Range("E19").Value = "  -3.99%  "

# Row 20
$ws.Range("D20").Value = "26.357.08"
$ws.Range("E20").Value = "  -3.00%  "

# Row 21
$ws.Range("D21").Value = "2.083.10"
$ws.Range("E21").Value = "  -1.83%  "

# Row 22
$ws.Range($scratch).Formula = "=""4.524"""
$ws.Range($scratch).Copy()
$ws.Range("D22").PasteSpecial(-4163)  # xlPasteValues
$ws.Range($scratch).Clear()
$ws.Range("E22").Value = "  -5.35%  "

# Row 23
$ws.Range($scratch).Formula = "=""9.391"""
$ws.Range($scratch).Copy()
$ws.Range("D23").PasteSpecial(-4163)  # xlPasteValues
$ws.Range($scratch).Clear()
$ws.Range("E23").Value = "  -7.53%  "

# Row 24
$ws.Range($scratch).Formula = "=""5.869"""
$ws.Range($scratch).Copy()
$ws.Range("D24").PasteSpecial(-4163)  # xlPasteValues
$ws.Range($scratch).Clear()
$ws.Range("E24").Value = "  -6.12%  "

# Row 25
$ws.Range($scratch).Formula = "=""2.279"""
$ws.Range($scratch).Copy()
$ws.Range("D25").PasteSpecial(-4163)  # xlPasteValues
$ws.Range($scratch).Clear()
$ws.Range("E25").Value = "  -5.38%  "

# Row 26
$ws.Range($scratch).Formula = "=""145.02"""
$ws.Range($scratch).Copy()
$ws.Range("D26").PasteSpecial(-4163)  # xlPasteValues
$ws.Range($scratch).Clear()
$ws.Range("E26").Value = "  -1.40%  "

# Row 27
$ws.Range($scratch).Formula = "=""1.643"""
$ws.Range($scratch).Copy()
$ws.Range("D27").PasteSpecial(-4163)  # xlPasteValues
$ws.Range($scratch).Clear()
$ws.Range("E27").Value = "  -1.97%  "

# Row 28
$ws.Range("E28").Value = "  -4.20%  "

# Row 29
$ws.Range($scratch).Formula = "=""110.73"""
$ws.Range($scratch).Copy()
$ws.Range("D29").PasteSpecial(-4163)  # xlPasteValues
$ws.Range($scratch).Clear()
$ws.Range("E29").Value = "  -2.27%  "

# Row 30
$ws.Range($scratch).Formula = "=""4.171"""
$ws.Range($scratch).Copy()
$ws.Range("D30").PasteSpecial(-4163)  # xlPasteValues
$ws.Range($scratch).Clear()
$ws.Range("E30").Value = "  -5.76%  "

# Row 31
$ws.Range($scratch).Formula = "=""4.106"""
$ws.Range($scratch).Copy()
$ws.Range("D31").PasteSpecial(-4163)  # xlPasteValues
$ws.Range($scratch).Clear()
$ws.Range("E31").Value = "  -6.50%  "

# Row 32
$ws.Range($scratch).Formula = "=""0.08700"""
$ws.Range($scratch).Copy()
$ws.Range("D32").PasteSpecial(-4163)  # xlPasteValues
$ws.Range($scratch).Clear()
$ws.Range("E32").Value = "  -2.86%  "

# Row 33
$ws.Range($scratch).Formula = "=""0.04820"""
$ws.Range($scratch).Copy()
$ws.Range("D33").PasteSpecial(-4163)  # xlPasteValues
$ws.Range($scratch).Clear()
$ws.Range("E33").Value = "  -2.33%  "

# Row 34
$ws.Range($scratch).Formula = "=""1.122"""
$ws.Range($scratch).Copy()
$ws.Range("D34").PasteSpecial(-4163)  # xlPasteValues
$ws.Range($scratch).Clear()
$ws.Range("E34").Value = "  -5.32%  "

# Row 35
$ws.Range($scratch).Formula = "=""2.845"""
$ws.Range($scratch).Copy()
$ws.Range("D35").PasteSpecial(-4163)  # xlPasteValues
$ws.Range($scratch).Clear()
$ws.Range("E35").Value = "  -2.20%  "

# Row 36
$ws.Range($scratch).Formula = "=""0.6810"""
$ws.Range($scratch).Copy()
$ws.Range("D36").PasteSpecial(-4163)  # xlPasteValues
$ws.Range($scratch).Clear()
$ws.Range("E36").Value = "  -9.15%  "

# Row 37
$ws.Range($scratch).Formula = "=""3.085"""
$ws.Range($scratch).Copy()
$ws.Range("D37").PasteSpecial(-4163)  # xlPasteValues
$ws.Range($scratch).Clear()
$ws.Range("E37").Value = "  -6.89%  "

# Row 38
$ws.Range($scratch).Formula = "=""0.01771"""
$ws.Range($scratch).Copy()
$ws.Range("D38").PasteSpecial(-4163)  # xlPasteValues
$ws.Range($scratch).Clear()
$ws.Range("E38").Value = "  -6.06%  "

# Row 39
$ws.Range($scratch).Formula = "=""2.200"""
$ws.Range($scratch).Copy()
$ws.Range("D39").PasteSpecial(-4163)  # xlPasteValues
$ws.Range($scratch).Clear()
$ws.Range("E39").Value = "  -9.12%  "

# Row 40
$ws.Range($scratch).Formula = "=""0.4882"""
$ws.Range($scratch).Copy()
$ws.Range("D40").PasteSpecial(-4163)  # xlPasteValues
$ws.Range($scratch).Clear()
$ws.Range("E40").Value = "  -8.67%  "

# Row 41
$ws.Range($scratch).Formula = "=""112.01"""
$ws.Range($scratch).Copy()
$ws.Range("D41").PasteSpecial(-4163)  # xlPasteValues
$ws.Range($scratch).Clear()
$ws.Range("E41").Value = "  -4.39%  "

# Row 42
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range($scratch).Formula = "=""6.119"""
$ws.Range($scratch).Copy()
$ws.Range("D42").PasteSpecial(-4163)  # xlPasteValues
$ws.Range($scratch).Clear()
$ws.Range("E42").Value = "  -3.47%  "

# Row 43
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range($scratch).Formula = "=""0.8839"""
$ws.Range($scratch).Copy()
$ws.Range("D43").PasteSpecial(-4163)  # xlPasteValues
$ws.Range($scratch).Clear()
$ws.Range("E43").Value = "  -9.13%  "

# Row 44
$ws.Range($scratch).Formula = "=""1.000"""
$ws.Range($scratch).Copy()
$ws.Range("D44").PasteSpecial(-4163)  # xlPasteValues
$ws.Range($scratch).Clear()
$ws.Range("E44").Value = "  +0.11%  "

# Row 45
$ws.Range($scratch).Formula = "=""7.658"""
$ws.Range($scratch).Copy()
$ws.Range("D45").PasteSpecial(-4163)  # xlPasteValues
$ws.Range($scratch).Clear()
$ws.Range("E45").Value = "  -7.01%  "

# Row 46
$ws.Range($scratch).Formula = "=""0.4180"""
$ws.Range($scratch).Copy()
$ws.Range("D46").PasteSpecial(-4163)  # xlPasteValues
$ws.Range($scratch).Clear()
$ws.Range("E46").Value = "  -9.62%  "

# Row 47
$ws.Range($scratch).Formula = "=""0.1250"""
$ws.Range($scratch).Copy()
$ws.Range("D47").PasteSpecial(-4163)  # xlPasteValues
$ws.Range($scratch).Clear()
$ws.Range("E47").Value = "  -8.91%  "

# Row 48
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range($scratch).Formula = "=""0.05872"""
$ws.Range($scratch).Copy()
$ws.Range("D48").PasteSpecial(-4163)  # xlPasteValues
$ws.Range($scratch).Clear()
$ws.Range("E48").Value = "  -1.44%  "

# Row 49
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range($scratch).Formula = "=""9.004"""
$ws.Range($scratch).Copy()
$ws.Range("D49").PasteSpecial(-4163)  # xlPasteValues
$ws.Range($scratch).Clear()
$ws.Range("E49").Value = "  -4.82%  "

# Row 50
$ws.Range($scratch).Formula = "=""35.20"""
$ws.Range($scratch).Copy()
$ws.Range("D50").PasteSpecial(-4163)  # xlPasteValues
$ws.Range($scratch).Clear()
$ws.Range("E50").Value = "  -4.37%  "

# Row 51
$ws.Range($scratch).Formula = "=""59.08"""
$ws.Range($scratch).Copy()
$ws.Range("D51").PasteSpecial(-4163)  # xlPasteValues
$ws.Range($scratch).Clear()
$ws.Range("E51").Value = "  -4.61%"
